$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 13: "ErrorHandler.reportError" -> "errorHandler.reportError"
# ---------------------------------------------------------------------------
$s13 = $p.Slides.Item(13)
$tr13 = $s13.Shapes.Item(2).TextFrame.TextRange
$found = $tr13.Find("ErrorHandler.reportError")
$found.Text = "errorHandler.reportError"

# ---------------------------------------------------------------------------
# Slide 15: "        ErrorHandler. " + "reportError"
#           -> "        " + "errorHandler.reportError"
# ---------------------------------------------------------------------------
$s15 = $p.Slides.Item(15)
$tf15 = $s15.Shapes.Item(2).TextFrame
$found = $tf15.TextRange.Find("        ErrorHandler. ")
$found.Text = "        "
$found = $tf15.TextRange.Find("reportError")
$found.Text = "errorHandler.reportError"

# ---------------------------------------------------------------------------
# Slide 18: "ErrorHandler.reportError" -> "errorHandler.reportError"
# ---------------------------------------------------------------------------
$s18 = $p.Slides.Item(18)
$tr18 = $s18.Shapes.Item(2).TextFrame.TextRange
$found = $tr18.Find("ErrorHandler.reportError")
$found.Text = "errorHandler.reportError"

# ---------------------------------------------------------------------------
# Slide 23: "ErrorHandler.reportError" -> "errorHandler.reportError"
# ---------------------------------------------------------------------------
$s23 = $p.Slides.Item(23)
$tr23 = $s23.Shapes.Item(2).TextFrame.TextRange
$found = $tr23.Find("ErrorHandler.reportError")
$found.Text = "errorHandler.reportError"

# ---------------------------------------------------------------------------
# Slide 7: title "Object " -> "Class "; body text clean-up
# ---------------------------------------------------------------------------
$s7 = $p.Slides.Item(7)

$titleTr = $s7.Shapes.Item(1).TextFrame.TextRange
$found = $titleTr.Find("Object ")
$found.Text = "Class "

$bodyTf = $s7.Shapes.Item(2).TextFrame
# Remove the whole "Kotlin object (not a class)." paragraph (first paragraph).
$para = $bodyTf.TextRange.Paragraphs(1, 1)
$para.Delete()

# Clear the (now 3rd) paragraph's text, leaving an empty paragraph in place of
# "Implements the singleton pattern; i.e., there is only one instance of ErrorHandler."
$para = $bodyTf.TextRange.Paragraphs(3, 1)
$para.Text = ""

# ---------------------------------------------------------------------------
# Slide 8: title "Two Key Methods in Object " -> "Two Key Methods in Class "
# ---------------------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$titleTr8 = $s8.Shapes.Item(1).TextFrame.TextRange
$found = $titleTr8.Find("Two Key Methods in Object ")
$found.Text = "Two Key Methods in Class "
